$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 85.31579000000001
$ws.Range("I11").Value = 85.31579000000001
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 85.31579000000001
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 54.68420999999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 129.4
$ws.Range("I33").Value = 136.28572
$ws.Range("J33").Value = 113.333336
$ws.Range("K33").Value = 136.28572
$ws.Range("L33").Value = 113.333336
$ws.Range("M33").Value = 92.71428
$ws.Range("N33").Value = -571.333336

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 315.8889
$ws.Range("I55").Value = 324.16666
$ws.Range("J55").Value = 299.33334
$ws.Range("K55").Value = 324.16666
$ws.Range("L55").Value = 299.33334
$ws.Range("M55").Value = -110.16666
$ws.Range("N55").Value = -727.33334

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1070.4375
$ws.Range("I98").Value = 995.13336
$ws.Range("J98").Value = 2200
$ws.Range("K98").Value = 995.13336
$ws.Range("L98").Value = 2200
$ws.Range("M98").Value = 502.86664
$ws.Range("N98").Value = -5196

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 3049.625
$ws.Range("I111").Value = 1200
$ws.Range("J111").Value = 3666.1667
$ws.Range("K111").Value = 3600
$ws.Range("L111").Value = 10998.5001
$ws.Range("M111").Value = -533
$ws.Range("N111").Value = -17132.5001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1070.4375
$ws.Range("I122").Value = 995.13336
$ws.Range("J122").Value = 2200
$ws.Range("K122").Value = 2985.40008
$ws.Range("L122").Value = 6600
$ws.Range("M122").Value = -535.4000800000003
$ws.Range("N122").Value = -11500

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2009.4634
$ws.Range("I138").Value = 1373.0555
$ws.Range("J138").Value = 2507.5217
$ws.Range("K138").Value = 4119.166499999999
$ws.Range("L138").Value = 7522.5651
$ws.Range("M138").Value = 1020.833500000001
$ws.Range("N138").Value = -17802.5651

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3797.6667
$ws.Range("I122").Value = 3797.6667
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 11393.0001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -8943.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2018.3438
$ws.Range("I132").Value = 2043.1852
$ws.Range("J132").Value = 1884.2
$ws.Range("K132").Value = 6129.5556
$ws.Range("L132").Value = 5652.6
$ws.Range("M132").Value = -3599.5556
$ws.Range("N132").Value = -10712.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 10484.333
$ws.Range("I105").Value = 8581.4
$ws.Range("J105").Value = 19999
$ws.Range("K105").Value = 8581.4
$ws.Range("L105").Value = 19999
$ws.Range("M105").Value = -6834.4
$ws.Range("N105").Value = -23493

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H88").Value = 25351.5
$ws.Range("I88").Value = 26311
$ws.Range("J88").Value = 25159.6
$ws.Range("K88").Value = 26311
$ws.Range("L88").Value = 25159.6
$ws.Range("M88").Value = -25905
$ws.Range("N88").Value = -25971.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H91").Value = 25351.5
$ws.Range("I91").Value = 26311
$ws.Range("J91").Value = 25159.6
$ws.Range("K91").Value = 26311
$ws.Range("L91").Value = 25159.6
$ws.Range("M91").Value = -24907
$ws.Range("N91").Value = -27967.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1265.125
$ws.Range("I107").Value = 945.8889
$ws.Range("J107").Value = 1456.6666
$ws.Range("K107").Value = 945.8889
$ws.Range("L107").Value = 1456.6666
$ws.Range("M107").Value = 974.1111
$ws.Range("N107").Value = -5296.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3199.6
$ws.Range("I134").Value = 2750
$ws.Range("J134").Value = 3499.3333
$ws.Range("K134").Value = 8250
$ws.Range("L134").Value = 10497.9999
$ws.Range("M134").Value = -5715
$ws.Range("N134").Value = -15567.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 582
$ws.Range("I18").Value = 582
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 1746
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -1577

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 79211
$ws.Range("I23").Value = 2166.3333
$ws.Range("J23").Value = 102324.4
$ws.Range("K23").Value = 6498.999899999999
$ws.Range("L23").Value = 306973.2
$ws.Range("M23").Value = -6263.999899999999
$ws.Range("N23").Value = -307443.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 12233.182
$ws.Range("I56").Value = 12233.182
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 12233.182
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -11703.182

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 1000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 1000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 3000
$ws.Range("N62").Value = -4372

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H65").Value = 1000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 1000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 9000
$ws.Range("N65").Value = -15864

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11831.333
$ws.Range("I70").Value = 10828.5
$ws.Range("J70").Value = 12834.167
$ws.Range("K70").Value = 10828.5
$ws.Range("L70").Value = 12834.167
$ws.Range("M70").Value = -10558.5
$ws.Range("N70").Value = -13374.167

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 11831.333
$ws.Range("I73").Value = 10828.5
$ws.Range("J73").Value = 12834.167
$ws.Range("K73").Value = 10828.5
$ws.Range("L73").Value = 12834.167
$ws.Range("M73").Value = -9892.5
$ws.Range("N73").Value = -14706.167

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2965.926
$ws.Range("I80").Value = 3002.8462
$ws.Range("J80").Value = 2006
$ws.Range("K80").Value = 3002.8462
$ws.Range("L80").Value = 2006
$ws.Range("M80").Value = -2004.8462
$ws.Range("N80").Value = -4002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2965.926
$ws.Range("I83").Value = 3002.8462
$ws.Range("J83").Value = 2006
$ws.Range("K83").Value = 15014.231
$ws.Range("L83").Value = 10030
$ws.Range("M83").Value = -10022.231
$ws.Range("N83").Value = -20014

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4899.933
$ws.Range("I126").Value = 4170.5
$ws.Range("J126").Value = 5165.1816
$ws.Range("K126").Value = 12511.5
$ws.Range("L126").Value = 15495.5448
$ws.Range("M126").Value = -10041.5
$ws.Range("N126").Value = -20435.5448

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 13199
$ws.Range("I16").Value = 8499
$ws.Range("J16").Value = 17899
$ws.Range("K16").Value = 8499
$ws.Range("L16").Value = 17899
$ws.Range("M16").Value = -8329
$ws.Range("N16").Value = -18239

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 101630.5
$ws.Range("I20").Value = 112833.89
$ws.Range("J20").Value = 800
$ws.Range("K20").Value = 112833.89
$ws.Range("L20").Value = 800
$ws.Range("M20").Value = -112607.89
$ws.Range("N20").Value = -1252

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1496.5385
$ws.Range("I22").Value = 1507
$ws.Range("J22").Value = 1491.8889
$ws.Range("K22").Value = 1507
$ws.Range("L22").Value = 1491.8889
$ws.Range("M22").Value = -1212
$ws.Range("N22").Value = -2081.8889

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 1496.5385
$ws.Range("I27").Value = 1507
$ws.Range("J27").Value = 1491.8889
$ws.Range("K27").Value = 1507
$ws.Range("L27").Value = 1491.8889
$ws.Range("M27").Value = -1400
$ws.Range("N27").Value = -1705.8889

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2183.111
$ws.Range("I61").Value = 1521.1428
$ws.Range("J61").Value = 4500
$ws.Range("K61").Value = 1521.1428
$ws.Range("L61").Value = 4500
$ws.Range("M61").Value = -1319.1428
$ws.Range("N61").Value = -4904

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 2183.111
$ws.Range("I113").Value = 1521.1428
$ws.Range("J113").Value = 4500
$ws.Range("K113").Value = 1521.1428
$ws.Range("L113").Value = 4500
$ws.Range("M113").Value = 648.8571999999999
$ws.Range("N113").Value = -8840

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 169134.83
$ws.Range("I132").Value = 169134.83
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 507404.49
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -504874.49

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2000000
$ws.Range("I2").Value = 2000000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 2000000
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1999888
$ws.Range("N2").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8308.333000000001
$ws.Range("I81").Value = 2450
$ws.Range("J81").Value = 14166.667
$ws.Range("K81").Value = 4900
$ws.Range("L81").Value = 28333.334
$ws.Range("M81").Value = -3839
$ws.Range("N81").Value = -30455.334

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 8308.333000000001
$ws.Range("I84").Value = 2450
$ws.Range("J84").Value = 14166.667
$ws.Range("K84").Value = 24500
$ws.Range("L84").Value = 141666.67
$ws.Range("M84").Value = -19196
$ws.Range("N84").Value = -152274.67

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4250
$ws.Range("I132").Value = 4000
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 12000
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -9470
$ws.Range("N132").Value = -20060
